$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "화성"
$ws.Range("C10").Value = "전곡항"
$ws.Range("D10").Value = "명성호"
$ws.Range("E10").Value = "http://xn--hq1b31ko5fzpfdsxrtb.com/index.php?mid=bk"
$ws.Range("F10").ClearContents()

# Row 11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "인천"
$ws.Range("C11").Value = "영흥항"
$ws.Range("D11").Value = "팀만수호"
$ws.Range("E11").Value = "https://teammansu.kr/index.php?mid=bk"
$ws.Range("F11").ClearContents()

# Row 12
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "군산"
$ws.Range("C12").Value = "비응항"
$ws.Range("D12").Value = "샤크호 선단"
$ws.Range("E12").Value = "https://sharkho.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F12").ClearContents()

# Row 13
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "화성"
$ws.Range("C13").Value = "전곡항"
$ws.Range("D13").Value = "신명호 선단"
$ws.Range("E13").Value = "http://www.shinmyungho.com/index.php?mid=bk"
$ws.Range("F13").ClearContents()

# Row 14
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "격포"
$ws.Range("C14").Value = "격포항"
$ws.Range("D14").Value = "변산레저낚시"
$ws.Range("E14").Value = "https://banak24.com/index.php?mid=bk"
$ws.Range("F14").Value = "카이저3호 문어만 이용 좋음"

# Row 15
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "군산"
$ws.Range("C15").Value = "비응항"
$ws.Range("D15").Value = "엘리스호"
$ws.Range("E15").Value = "https://alice.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F15").ClearContents()

# Row 16
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "보령"
$ws.Range("C16").Value = "홍원항"
$ws.Range("D16").Value = "뉴해양호"
$ws.Range("E16").Value = "https://newhaeyang.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F16").ClearContents()

# Row 17
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "보령"
$ws.Range("C17").Value = "홍원항"
$ws.Range("D17").Value = "조커호"
$ws.Range("E17").Value = "http://seasidefishing.kr/index.php?mid=bk"
$ws.Range("F17").Value = "BDJ 추천 배, 선장 젊고, 열정 많음
"

# Row 18
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "보령"
$ws.Range("C18").Value = "대천항"
$ws.Range("D18").Value = "팀루피호"
$ws.Range("E18").Value = "https://masterfishing.kr/index.php?mid=bk"
$ws.Range("F18").ClearContents()

# Row 19
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "보령"
$ws.Range("C19").Value = "구매항"
$ws.Range("D19").Value = "악바리호"
$ws.Range("E19").Value = "https://akbari.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F19").ClearContents()

# Row 20
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "보령"
$ws.Range("C20").Value = "오천항"
$ws.Range("D20").Value = "범블비호"
$ws.Range("E20").Value = "http://xn--xk3bm1aee249g.com/index.php?mid=bk"
$ws.Range("F20").ClearContents()

# Row 21
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "보령"
$ws.Range("C21").Value = "오천항"
$ws.Range("D21").Value = "루키호"
$ws.Range("E21").Value = "http://www.yamujinfishing.com/index.php?mid=bk"
$ws.Range("F21").Value = "오천에 유명 배"

# Row 22
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "보령"
$ws.Range("C22").Value = "오천항"
$ws.Range("D22").Value = "프랜드피싱"
$ws.Range("E22").Value = "http://www.friendho.com/index.php?mid=bk"
$ws.Range("F22").Value = "멘구전문 갑,쭈 잘함"

# Row 23
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "태안"
$ws.Range("C23").Value = "마검포항"
$ws.Range("D23").Value = "가가호"
$ws.Range("E23").Value = "https://gagaho.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F23").Value = "내꼬야,유튜브가 주로 이용, 갑잘함"

# Row 24
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "태안"
$ws.Range("C24").Value = "마검포항"
$ws.Range("D24").Value = "뉴정원호"
$ws.Range("E24").Value = "https://www.jungwonho.com/index.php?mid=bk"
$ws.Range("F24").ClearContents()

# Row 25
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "태안"
$ws.Range("C25").Value = "마검포항"
$ws.Range("D25").Value = "만선호"
$ws.Range("E25").Value = "https://mansunho.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F25").ClearContents()

# Row 26
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "태안"
$ws.Range("C26").Value = "영목항"
$ws.Range("D26").Value = "카라호"
$ws.Range("E26").Value = "https://karaho.kr/index.php?mid=bk"
$ws.Range("F26").ClearContents()

# Row 27
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "태안"
$ws.Range("C27").Value = "영목항"
$ws.Range("D27").Value = "청광호"
$ws.Range("E27").Value = "http://www.chungkwangho.net/index.php?mid=bk"
$ws.Range("F27").Value = "우럭귀신 자주 이용하는 배"

# Row 28
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "당진"
$ws.Range("C28").Value = "장고항"
$ws.Range("D28").Value = "디스커버리호"
$ws.Range("E28").Value = "http://www.discoveryho.net/index.php?mid=bk"
$ws.Range("F28").Value = "김조사, 송조사 추천, 꼬리쪽 잘함"

# Row 29
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "안산"
$ws.Range("C29").Value = "오이도항"
$ws.Range("D29").Value = "포세이돈호"
$ws.Range("E29").Value = "https://poseidon.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F29").ClearContents()

# Row 30
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "인천"
$ws.Range("C30").Value = "영흥항"
$ws.Range("D30").Value = "금강7호"
$ws.Range("E30").Value = "http://www.kumkangho.co.kr/index.php?mid=bk"
$ws.Range("F30").ClearContents()

# Row 31
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "여수"
$ws.Range("C31").Value = "국동항"
$ws.Range("D31").Value = "트윈스타"
$ws.Range("E31").Value = "https://twin.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F31").ClearContents()

# Row 32
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "여수"
$ws.Range("C32").Value = "국동항"
$ws.Range("D32").Value = "만진스타호"
$ws.Range("E32").Value = "https://manjinstar.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F32").ClearContents()

# Row 33
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "여수"
$ws.Range("C33").Value = "국동항"
$ws.Range("D33").Value = "오션스타호"
$ws.Range("E33").Value = "https://ysoceanstar.sunsang24.com/ship/schedule_fleet"
$ws.Range("F33").ClearContents()

# Row 34
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "여수"
$ws.Range("C34").Value = "종포항"
$ws.Range("D34").Value = "빅보스호"
$ws.Range("E34").Value = "https://big.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F34").ClearContents()

# Row 35
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "여수"
$ws.Range("C35").Value = "국동항"
$ws.Range("D35").Value = "써니호"
$ws.Range("E35").Value = "https://sunny.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F35").ClearContents()

# Row 36
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "고흥"
$ws.Range("C36").Value = "녹동방파제"
$ws.Range("D36").Value = "빅원호"
$ws.Range("E36").Value = "https://bigone.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F36").ClearContents()

# Row 37
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "고흥"
$ws.Range("C37").Value = "녹동방파제"
$ws.Range("D37").Value = "에이스호"
$ws.Range("E37").Value = "https://aceho.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F37").Value = "신조선, 선장 배질 잘함,깨끗함"

# Row 38
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "여수"
$ws.Range("C38").Value = "신추항"
$ws.Range("D38").Value = "그린나래"
$ws.Range("E38").Value = "https://green.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F38").ClearContents()

# Row 39
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "서산"
$ws.Range("C39").Value = "삼길포항"
$ws.Range("D39").Value = "만석낚시 선단"
$ws.Range("E39").Value = "http://www.mscufishing.com/index.php?mid=bk"
$ws.Range("F39").ClearContents()

# Row 40
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "여수"
$ws.Range("C40").Value = "돌산항"
$ws.Range("D40").Value = "영심이호"
$ws.Range("E40").Value = "https://0simi.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F40").Value = "전투 낚시 전문, 배는 작지만, 선장 마인드 좋음"

# Row 41
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "인천"
$ws.Range("C41").Value = "영흥항"
$ws.Range("D41").Value = "만수피싱"
$ws.Range("E41").Value = "https://teammansu.kr/index.php?mid=bk"
$ws.Range("F41").ClearContents()

# Row 42
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "보령"
$ws.Range("C42").Value = "남당항"
$ws.Range("D42").Value = "은가비호"
$ws.Range("E42").Value = "https://eungabi.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F42").ClearContents()

# Row 43
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "서산"
$ws.Range("C43").Value = "삼길포항"
$ws.Range("D43").Value = "헤르메스호"
$ws.Range("E43").Value = "http://hermes.thefishing.kr/index.php?mid=bk"
$ws.Range("F43").ClearContents()

# Row 44
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "평택"
$ws.Range("C44").Value = "평택항"
$ws.Range("D44").Value = "수복호"
$ws.Range("E44").Value = "https://subokho.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F44").ClearContents()

# Row 45
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "태안"
$ws.Range("C45").Value = "모항항"
$ws.Range("D45").Value = "라온피싱"
$ws.Range("E45").Value = "http://raonfishing.com/index.php?mid=bk"
$ws.Range("F45").ClearContents()

# Row 46
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "고흥"
$ws.Range("C46").Value = "녹동방파제"
$ws.Range("D46").Value = "몬스터호"
$ws.Range("E46").Value = "https://mon.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F46").ClearContents()

# Row 47
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "고흥"
$ws.Range("C47").Value = "녹동방파제"
$ws.Range("D47").Value = "여명호"
$ws.Range("E47").Value = "https://ym.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F47").ClearContents()

# Row 48
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "여수"
$ws.Range("C48").Value = "돌산항"
$ws.Range("D48").Value = "청홍낚시"
$ws.Range("E48").Value = "https://chf.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F48").ClearContents()

# Row 49
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "여수"
$ws.Range("C49").Value = "국동항"
$ws.Range("D49").Value = "블랙펄호"
$ws.Range("E49").Value = "https://ysblackpearl.sunsang24.com/ship/schedule_fleet/"
$ws.Range("F49").Value = "멘구 추천 문어 선사"
